$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E24").Value = 29
$ws.Range("E26").Value = 30

$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(4)
$series.Values = "=Sheet1!`$E`$2:`$E`$26"
Write-Output $series.Values
